# Fix GeneratorsOfNode sheet: include both "Wind offshore grounded" and
# "Wind offshore floating" rows for every offshore wind area (UK + Norwegian
# areas). Previously rows 738-773 had a mixed/partial listing; now rows
# 738-768 list all 31 areas as grounded and rows 769-799 list the same 31
# areas as floating.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneratorsOfNode")

$data = @(
    @(738, 'Moray Firth', 'Wind offshore grounded'),
    @(739, 'Firth of Forth', 'Wind offshore grounded'),
    @(740, 'Dogger Bank', 'Wind offshore grounded'),
    @(741, 'Hornsea', 'Wind offshore grounded'),
    @(742, 'Outer Dowsing', 'Wind offshore grounded'),
    @(743, 'Norfolk', 'Wind offshore grounded'),
    @(744, 'East Anglia', 'Wind offshore grounded'),
    @(745, 'Borssele', 'Wind offshore grounded'),
    @(746, 'Hollandsee Kust', 'Wind offshore grounded'),
    @(747, 'Helgoländer Bucht', 'Wind offshore grounded'),
    @(748, 'Nordsøen', 'Wind offshore grounded'),
    @(749, 'Nordvest A', 'Wind offshore grounded'),
    @(750, 'Nordvest C', 'Wind offshore grounded'),
    @(751, 'Vestavind A', 'Wind offshore grounded'),
    @(752, 'Sønnavind A', 'Wind offshore grounded'),
    @(753, 'Sørvest C', 'Wind offshore grounded'),
    @(754, 'Nordvest B', 'Wind offshore grounded'),
    @(755, 'Vestavind F', 'Wind offshore grounded'),
    @(756, 'Vestavind B', 'Wind offshore grounded'),
    @(757, 'Vestavind C', 'Wind offshore grounded'),
    @(758, 'Vestavind D', 'Wind offshore grounded'),
    @(759, 'Sørvest F', 'Wind offshore grounded'),
    @(760, 'Sørvest B', 'Wind offshore grounded'),
    @(761, 'Nordavind B', 'Wind offshore grounded'),
    @(762, 'Nordavind A', 'Wind offshore grounded'),
    @(763, 'Nordavind D', 'Wind offshore grounded'),
    @(764, 'Nordavind C', 'Wind offshore grounded'),
    @(765, 'Vestavind E', 'Wind offshore grounded'),
    @(766, 'Sørvest E', 'Wind offshore grounded'),
    @(767, 'Sørvest A', 'Wind offshore grounded'),
    @(768, 'Sørvest D', 'Wind offshore grounded'),
    @(769, 'Moray Firth', 'Wind offshore floating'),
    @(770, 'Firth of Forth', 'Wind offshore floating'),
    @(771, 'Dogger Bank', 'Wind offshore floating'),
    @(772, 'Hornsea', 'Wind offshore floating'),
    @(773, 'Outer Dowsing', 'Wind offshore floating'),
    @(774, 'Norfolk', 'Wind offshore floating'),
    @(775, 'East Anglia', 'Wind offshore floating'),
    @(776, 'Borssele', 'Wind offshore floating'),
    @(777, 'Hollandsee Kust', 'Wind offshore floating'),
    @(778, 'Helgoländer Bucht', 'Wind offshore floating'),
    @(779, 'Nordsøen', 'Wind offshore floating'),
    @(780, 'Nordvest A', 'Wind offshore floating'),
    @(781, 'Nordvest C', 'Wind offshore floating'),
    @(782, 'Vestavind A', 'Wind offshore floating'),
    @(783, 'Sønnavind A', 'Wind offshore floating'),
    @(784, 'Sørvest C', 'Wind offshore floating'),
    @(785, 'Nordvest B', 'Wind offshore floating'),
    @(786, 'Vestavind F', 'Wind offshore floating'),
    @(787, 'Vestavind B', 'Wind offshore floating'),
    @(788, 'Vestavind C', 'Wind offshore floating'),
    @(789, 'Vestavind D', 'Wind offshore floating'),
    @(790, 'Sørvest F', 'Wind offshore floating'),
    @(791, 'Sørvest B', 'Wind offshore floating'),
    @(792, 'Nordavind B', 'Wind offshore floating'),
    @(793, 'Nordavind A', 'Wind offshore floating'),
    @(794, 'Nordavind D', 'Wind offshore floating'),
    @(795, 'Nordavind C', 'Wind offshore floating'),
    @(796, 'Vestavind E', 'Wind offshore floating'),
    @(797, 'Sørvest E', 'Wind offshore floating'),
    @(798, 'Sørvest A', 'Wind offshore floating'),
    @(799, 'Sørvest D', 'Wind offshore floating')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
